# Update Leave Card 6/29/2023 4:54 PM
#
# Fills in the PERIOD (date) column for the monthly rows that were added to
# the leave card table (Mar-2023 .. Mar-2024), and records the 1.25 EARNED
# credit for Mar/Apr/May 2023 (the three months that had actually accrued
# leave by the time of this save). The table's calculated columns (EARNED /
# BALANCE) recompute automatically from these inputs.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("CONVERTION")

# --- PERIOD dates for the newly-extended table rows (row 32 = Mar 2023) ---
$periods = @{
    32 = 44986  # 3/1/2023
    33 = 45017  # 4/1/2023
    34 = 45047  # 5/1/2023
    35 = 45078  # 6/1/2023
    36 = 45108  # 7/1/2023
    37 = 45139  # 8/1/2023
    38 = 45170  # 9/1/2023
    39 = 45200  # 10/1/2023
    40 = 45231  # 11/1/2023
    41 = 45261  # 12/1/2023
    42 = 45292  # 1/1/2024
    43 = 45323  # 2/1/2024
    44 = 45352  # 3/1/2024
}

foreach ($row in $periods.Keys) {
    $ws1.Cells.Item($row, 1).Value = $periods[$row]
}

# --- EARNED credits of 1.25 for Mar / Apr / May 2023 (rows 32-34) ---
$ws1.Cells.Item(32, 3).Value = 1.25
$ws1.Cells.Item(33, 3).Value = 1.25
$ws1.Cells.Item(34, 3).Value = 1.25

# --- Restore the view/selection state recorded in the saved workbook ---
$ws2.Select() | Out-Null
$ws2.Range("I20").Select() | Out-Null

$ws1.Select() | Out-Null
$ws1.Range("H25").Select() | Out-Null
